$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new row (39) date cell uses the same date style as the existing date column
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 2
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44445
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 100114007
$ws.Cells.Item(2, 7).Value = "Jengibre"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 600
$ws.Cells.Item(2, 11).Value = 13000
$ws.Cells.Item(2, 12).Value = 14000
$ws.Cells.Item(2, 13).Value = 13500
$ws.Cells.Item(2, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(2, 15).Value = "Perú"
$ws.Cells.Item(2, 16).Value = 1038
$ws.Cells.Item(2, 17).Value = 13
$ws.Cells.Item(2, 18).Value = "Hortaliza"

# Row 3
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44658
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = "Jengibre"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 400
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 1192
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# Row 4
$ws.Cells.Item(4, 1).Value = 8
$ws.Cells.Item(4, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(4, 3).Value = "Coquimbo"
$ws.Cells.Item(4, 4).Value = 44335
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 100114007
$ws.Cells.Item(4, 7).Value = "Jengibre"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 480
$ws.Cells.Item(4, 11).Value = 24500
$ws.Cells.Item(4, 12).Value = 25000
$ws.Cells.Item(4, 13).Value = 24750
$ws.Cells.Item(4, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(4, 15).Value = "Perú"
$ws.Cells.Item(4, 16).Value = 1904
$ws.Cells.Item(4, 17).Value = 13
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Row 5
$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44428
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 100114007
$ws.Cells.Item(5, 7).Value = "Jengibre"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 480
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14500
$ws.Cells.Item(5, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 1115
$ws.Cells.Item(5, 17).Value = 13
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# Row 6
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44442
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 100114007
$ws.Cells.Item(6, 7).Value = "Jengibre"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 460
$ws.Cells.Item(6, 11).Value = 14000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 13).Value = 14500
$ws.Cells.Item(6, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(6, 15).Value = "Perú"
$ws.Cells.Item(6, 16).Value = 1115
$ws.Cells.Item(6, 17).Value = 13
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Row 7
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44599
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100114007
$ws.Cells.Item(7, 7).Value = "Jengibre"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 15000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 15500
$ws.Cells.Item(7, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 1192
$ws.Cells.Item(7, 17).Value = 13
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Row 8
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44435
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100114007
$ws.Cells.Item(8, 7).Value = "Jengibre"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 480
$ws.Cells.Item(8, 11).Value = 13000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 13500
$ws.Cells.Item(8, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 1038
$ws.Cells.Item(8, 17).Value = 13
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44309
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 100114007
$ws.Cells.Item(9, 7).Value = "Jengibre"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 26000
$ws.Cells.Item(9, 12).Value = 27000
$ws.Cells.Item(9, 13).Value = 26500
$ws.Cells.Item(9, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 2038
$ws.Cells.Item(9, 17).Value = 13
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44400
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 600
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 16000
$ws.Cells.Item(10, 13).Value = 15500
$ws.Cells.Item(10, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 1192
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# Row 11
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44344
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = "Jengibre"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 400
$ws.Cells.Item(11, 11).Value = 18500
$ws.Cells.Item(11, 12).Value = 19000
$ws.Cells.Item(11, 13).Value = 18750
$ws.Cells.Item(11, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 1442
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Row 12
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44484
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100114007
$ws.Cells.Item(12, 7).Value = "Jengibre"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 360
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14500
$ws.Cells.Item(12, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 1115
$ws.Cells.Item(12, 17).Value = 13
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# Row 13
$ws.Cells.Item(13, 1).Value = 8
$ws.Cells.Item(13, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44410
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100114007
$ws.Cells.Item(13, 7).Value = "Jengibre"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 600
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 1115
$ws.Cells.Item(13, 17).Value = 13
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Row 14
$ws.Cells.Item(14, 1).Value = 8
$ws.Cells.Item(14, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44333
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 100114007
$ws.Cells.Item(14, 7).Value = "Jengibre"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 440
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 25000
$ws.Cells.Item(14, 13).Value = 24500
$ws.Cells.Item(14, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 1885
$ws.Cells.Item(14, 17).Value = 13
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Row 15
$ws.Cells.Item(15, 1).Value = 8
$ws.Cells.Item(15, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44680
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = 100114007
$ws.Cells.Item(15, 7).Value = "Jengibre"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 13500
$ws.Cells.Item(15, 12).Value = 14000
$ws.Cells.Item(15, 13).Value = 13750
$ws.Cells.Item(15, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 1058
$ws.Cells.Item(15, 17).Value = 13
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Row 16
$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44323
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = "Jengibre"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 460
$ws.Cells.Item(16, 11).Value = 25000
$ws.Cells.Item(16, 12).Value = 26000
$ws.Cells.Item(16, 13).Value = 25500
$ws.Cells.Item(16, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1962
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Row 17
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44670
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100114007
$ws.Cells.Item(17, 7).Value = "Jengibre"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 480
$ws.Cells.Item(17, 11).Value = 14500
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 14750
$ws.Cells.Item(17, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 1135
$ws.Cells.Item(17, 17).Value = 13
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Row 18
$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44326
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 100114007
$ws.Cells.Item(18, 7).Value = "Jengibre"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 460
$ws.Cells.Item(18, 11).Value = 25000
$ws.Cells.Item(18, 12).Value = 26000
$ws.Cells.Item(18, 13).Value = 25500
$ws.Cells.Item(18, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 1962
$ws.Cells.Item(18, 17).Value = 13
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Row 19
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44687
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 100114007
$ws.Cells.Item(19, 7).Value = "Jengibre"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 440
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14500
$ws.Cells.Item(19, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 1115
$ws.Cells.Item(19, 17).Value = 13
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# Row 20
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44498
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 100114007
$ws.Cells.Item(20, 7).Value = "Jengibre"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 14500
$ws.Cells.Item(20, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 1115
$ws.Cells.Item(20, 17).Value = 13
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# Row 21
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44582
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 100114007
$ws.Cells.Item(21, 7).Value = "Jengibre"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 520
$ws.Cells.Item(21, 11).Value = 15000
$ws.Cells.Item(21, 12).Value = 16000
$ws.Cells.Item(21, 13).Value = 15500
$ws.Cells.Item(21, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 1192
$ws.Cells.Item(21, 17).Value = 13
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# Row 22
$ws.Cells.Item(22, 1).Value = 8
$ws.Cells.Item(22, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44312
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = 100114007
$ws.Cells.Item(22, 7).Value = "Jengibre"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 26000
$ws.Cells.Item(22, 12).Value = 27000
$ws.Cells.Item(22, 13).Value = 26500
$ws.Cells.Item(22, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(22, 15).Value = "Perú"
$ws.Cells.Item(22, 16).Value = 2038
$ws.Cells.Item(22, 17).Value = 13
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# Row 23
$ws.Cells.Item(23, 1).Value = 8
$ws.Cells.Item(23, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value = 44418
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = 100114007
$ws.Cells.Item(23, 7).Value = "Jengibre"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 500
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 14500
$ws.Cells.Item(23, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(23, 15).Value = "Perú"
$ws.Cells.Item(23, 16).Value = 1115
$ws.Cells.Item(23, 17).Value = 13
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Row 24
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44412
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100114007
$ws.Cells.Item(24, 7).Value = "Jengibre"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 600
$ws.Cells.Item(24, 11).Value = 14000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 14500
$ws.Cells.Item(24, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(24, 15).Value = "Perú"
$ws.Cells.Item(24, 16).Value = 1115
$ws.Cells.Item(24, 17).Value = 13
$ws.Cells.Item(24, 18).Value = "Hortaliza"

# Row 25
$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44692
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = 100114007
$ws.Cells.Item(25, 7).Value = "Jengibre"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 400
$ws.Cells.Item(25, 11).Value = 14000
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14500
$ws.Cells.Item(25, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 1115
$ws.Cells.Item(25, 17).Value = 13
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Row 26
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44533
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100114007
$ws.Cells.Item(26, 7).Value = "Jengibre"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 520
$ws.Cells.Item(26, 11).Value = 17000
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = 17500
$ws.Cells.Item(26, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 1346
$ws.Cells.Item(26, 17).Value = 13
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Row 27
$ws.Cells.Item(27, 1).Value = 8
$ws.Cells.Item(27, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44414
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = 100114007
$ws.Cells.Item(27, 7).Value = "Jengibre"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 500
$ws.Cells.Item(27, 11).Value = 14000
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = 14500
$ws.Cells.Item(27, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 1115
$ws.Cells.Item(27, 17).Value = 13
$ws.Cells.Item(27, 18).Value = "Hortaliza"

# Row 28
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44596
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100114007
$ws.Cells.Item(28, 7).Value = "Jengibre"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 500
$ws.Cells.Item(28, 11).Value = 16000
$ws.Cells.Item(28, 12).Value = 17000
$ws.Cells.Item(28, 13).Value = 16500
$ws.Cells.Item(28, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 1269
$ws.Cells.Item(28, 17).Value = 13
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Row 29
$ws.Cells.Item(29, 1).Value = 8
$ws.Cells.Item(29, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(29, 3).Value = "Coquimbo"
$ws.Cells.Item(29, 4).Value = 44657
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = 100114007
$ws.Cells.Item(29, 7).Value = "Jengibre"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 460
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 16000
$ws.Cells.Item(29, 13).Value = 15500
$ws.Cells.Item(29, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 1192
$ws.Cells.Item(29, 17).Value = 13
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Row 30
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44426
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100114007
$ws.Cells.Item(30, 7).Value = "Jengibre"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 460
$ws.Cells.Item(30, 11).Value = 14000
$ws.Cells.Item(30, 12).Value = 15000
$ws.Cells.Item(30, 13).Value = 14500
$ws.Cells.Item(30, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(30, 15).Value = "Perú"
$ws.Cells.Item(30, 16).Value = 1115
$ws.Cells.Item(30, 17).Value = 13
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# Row 31
$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44575
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100114007
$ws.Cells.Item(31, 7).Value = "Jengibre"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 500
$ws.Cells.Item(31, 11).Value = 14000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 14500
$ws.Cells.Item(31, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(31, 15).Value = "Perú"
$ws.Cells.Item(31, 16).Value = 1115
$ws.Cells.Item(31, 17).Value = 13
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Row 32
$ws.Cells.Item(32, 1).Value = 8
$ws.Cells.Item(32, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44379
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100114007
$ws.Cells.Item(32, 7).Value = "Jengibre"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 600
$ws.Cells.Item(32, 11).Value = 17000
$ws.Cells.Item(32, 12).Value = 18000
$ws.Cells.Item(32, 13).Value = 17500
$ws.Cells.Item(32, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(32, 15).Value = "Perú"
$ws.Cells.Item(32, 16).Value = 1346
$ws.Cells.Item(32, 17).Value = 13
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Row 33
$ws.Cells.Item(33, 1).Value = 8
$ws.Cells.Item(33, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 44505
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 100114007
$ws.Cells.Item(33, 7).Value = "Jengibre"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 400
$ws.Cells.Item(33, 11).Value = 16000
$ws.Cells.Item(33, 12).Value = 17000
$ws.Cells.Item(33, 13).Value = 16500
$ws.Cells.Item(33, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 1269
$ws.Cells.Item(33, 17).Value = 13
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Row 34
$ws.Cells.Item(34, 1).Value = 8
$ws.Cells.Item(34, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44260
$ws.Cells.Item(34, 5).Value = 4
$ws.Cells.Item(34, 6).Value = 100114007
$ws.Cells.Item(34, 7).Value = "Jengibre"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 400
$ws.Cells.Item(34, 11).Value = 37000
$ws.Cells.Item(34, 12).Value = 38000
$ws.Cells.Item(34, 13).Value = 37500
$ws.Cells.Item(34, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(34, 15).Value = "Perú"
$ws.Cells.Item(34, 16).Value = 2885
$ws.Cells.Item(34, 17).Value = 13
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Row 35
$ws.Cells.Item(35, 1).Value = 8
$ws.Cells.Item(35, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(35, 3).Value = "Coquimbo"
$ws.Cells.Item(35, 4).Value = 44383
$ws.Cells.Item(35, 5).Value = 4
$ws.Cells.Item(35, 6).Value = 100114007
$ws.Cells.Item(35, 7).Value = "Jengibre"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 200
$ws.Cells.Item(35, 11).Value = 17000
$ws.Cells.Item(35, 12).Value = 18000
$ws.Cells.Item(35, 13).Value = 17500
$ws.Cells.Item(35, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(35, 15).Value = "Perú"
$ws.Cells.Item(35, 16).Value = 1346
$ws.Cells.Item(35, 17).Value = 13
$ws.Cells.Item(35, 18).Value = "Hortaliza"

# Row 36
$ws.Cells.Item(36, 1).Value = 8
$ws.Cells.Item(36, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 44644
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(36, 6).Value = 100114007
$ws.Cells.Item(36, 7).Value = "Jengibre"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 400
$ws.Cells.Item(36, 11).Value = 15000
$ws.Cells.Item(36, 12).Value = 16000
$ws.Cells.Item(36, 13).Value = 15500
$ws.Cells.Item(36, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(36, 15).Value = "Perú"
$ws.Cells.Item(36, 16).Value = 1192
$ws.Cells.Item(36, 17).Value = 13
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# Row 37
$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44631
$ws.Cells.Item(37, 5).Value = 4
$ws.Cells.Item(37, 6).Value = 100114007
$ws.Cells.Item(37, 7).Value = "Jengibre"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 16000
$ws.Cells.Item(37, 12).Value = 17000
$ws.Cells.Item(37, 13).Value = 16500
$ws.Cells.Item(37, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(37, 15).Value = "Perú"
$ws.Cells.Item(37, 16).Value = 1269
$ws.Cells.Item(37, 17).Value = 13
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Row 38
$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44365
$ws.Cells.Item(38, 5).Value = 4
$ws.Cells.Item(38, 6).Value = 100114007
$ws.Cells.Item(38, 7).Value = "Jengibre"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 500
$ws.Cells.Item(38, 11).Value = 19500
$ws.Cells.Item(38, 12).Value = 20000
$ws.Cells.Item(38, 13).Value = 19750
$ws.Cells.Item(38, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 1519
$ws.Cells.Item(38, 17).Value = 13
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# Row 39
$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44419
$ws.Cells.Item(39, 5).Value = 4
$ws.Cells.Item(39, 6).Value = 100114007
$ws.Cells.Item(39, 7).Value = "Jengibre"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 600
$ws.Cells.Item(39, 11).Value = 14000
$ws.Cells.Item(39, 12).Value = 15000
$ws.Cells.Item(39, 13).Value = 14500
$ws.Cells.Item(39, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(39, 15).Value = "Perú"
$ws.Cells.Item(39, 16).Value = 1115
$ws.Cells.Item(39, 17).Value = 13
$ws.Cells.Item(39, 18).Value = "Hortaliza"
